$xlShiftDown    = -4121
$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# The "Metadata" sheet lists one row per property of the StructureDefinition.
# It already has one "Contact" row (row 11, with row 10 being a duplicate of
# it). This commit adds all of the IG's authors as contacts, which means two
# more "Contact" / "No display for ContactDetail" rows need to be inserted
# right after the existing ones (i.e. as new rows 12 and 13), pushing every
# following property row down by two.
$ws.Range("A12:B13").Insert($xlShiftDown)

# Copy the formatting from the existing Contact row (11) onto the two new
# rows so they keep the same style as the rest of the table.
$ws.Range("A11:B11").Copy()
$ws.Range("A13:B13").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("A12").Value = "Contact"
$ws.Range("B12").Value = "No display for ContactDetail"
$ws.Range("A13").Value = "Contact"
$ws.Range("B13").Value = "No display for ContactDetail"

# Updated "Date" value, a side effect of regenerating the IG export.
$ws.Range("B8").Value = "2022-01-21T07:49:24+01:00"
